$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 202.625
$ws.Cells.Item(8, 9).Value = 231.42857
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = 694.28571
$ws.Cells.Item(8, 12).Value = 3
$ws.Cells.Item(8, 13).Value = -555.28571
$ws.Cells.Item(8, 14).Value = -281
$ws.Cells.Item(9, 8).Value = 136.5
$ws.Cells.Item(9, 9).Value = 188
$ws.Cells.Item(9, 11).Value = 188
$ws.Cells.Item(9, 13).Value = -19
$ws.Cells.Item(33, 8).Value = 164.22223
$ws.Cells.Item(33, 9).Value = 136.8
$ws.Cells.Item(33, 11).Value = 136.8
$ws.Cells.Item(33, 13).Value = 92.19999999999999
$ws.Cells.Item(116, 8).Value = 2201
$ws.Cells.Item(116, 9).Value = 1999.8572
$ws.Cells.Item(116, 11).Value = 1999.8572
$ws.Cells.Item(116, 13).Value = 1442.1428
$ws.Cells.Item(135, 8).Value = 969.3333
$ws.Cells.Item(135, 9).Value = 1003
$ws.Cells.Item(135, 11).Value = 9027
$ws.Cells.Item(135, 13).Value = -6492

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 18000
$ws.Cells.Item(23, 10).Value = 18000
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 14).Value = -18518
$ws.Cells.Item(114, 8).Value = 15000000
$ws.Cells.Item(114, 10).Value = 15000000
$ws.Cells.Item(114, 12).Value = 15000000
$ws.Cells.Item(114, 14).Value = -15008678
$ws.Cells.Item(132, 8).Value = 3794.926
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = $null

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 7099.5
$ws.Cells.Item(86, 9).Value = 5398.5
$ws.Cells.Item(86, 10).Value = 7950
$ws.Cells.Item(86, 11).Value = 5398.5
$ws.Cells.Item(86, 12).Value = 7950
$ws.Cells.Item(86, 13).Value = -4275.5
$ws.Cells.Item(86, 14).Value = -10196
$ws.Cells.Item(89, 8).Value = 7099.5
$ws.Cells.Item(89, 9).Value = 5398.5
$ws.Cells.Item(89, 10).Value = 7950
$ws.Cells.Item(89, 11).Value = 26992.5
$ws.Cells.Item(89, 12).Value = 39750
$ws.Cells.Item(89, 13).Value = -21376.5
$ws.Cells.Item(89, 14).Value = -50982
$ws.Cells.Item(107, 8).Value = 4425.0625
$ws.Cells.Item(107, 9).Value = 3065.6667
$ws.Cells.Item(107, 11).Value = 3065.6667
$ws.Cells.Item(107, 13).Value = -1145.6667
$ws.Cells.Item(134, 8).Value = 2284.5
$ws.Cells.Item(134, 9).Value = 2284.5
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 6853.5
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -4318.5
$ws.Cells.Item(134, 14).Value = $null

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(21, 8).Value = 35173.332
$ws.Cells.Item(21, 10).Value = 24573.572
$ws.Cells.Item(21, 12).Value = 24573.572
$ws.Cells.Item(21, 14).Value = -25043.572
$ws.Cells.Item(134, 8).Value = 2079.2222
$ws.Cells.Item(134, 9).Value = 2079.2222
$ws.Cells.Item(134, 11).Value = 6237.6666
$ws.Cells.Item(134, 13).Value = -3702.6666
$ws.Cells.Item(140, 8).Value = 56149
$ws.Cells.Item(140, 10).Value = 56149
$ws.Cells.Item(140, 12).Value = 56149
$ws.Cells.Item(140, 14).Value = -66509

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 78.22727
$ws.Cells.Item(2, 10).Value = 41.375
$ws.Cells.Item(2, 12).Value = 248.25
$ws.Cells.Item(2, 14).Value = -474.25
$ws.Cells.Item(4, 8).Value = 44519.934
$ws.Cells.Item(4, 10).Value = 75.2
$ws.Cells.Item(4, 12).Value = 225.6
$ws.Cells.Item(4, 14).Value = -449.6
$ws.Cells.Item(6, 8).Value = 58.444443
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(6, 12).Value = 180
$ws.Cells.Item(6, 14).Value = -406
$ws.Cells.Item(14, 8).Value = 740.8333
$ws.Cells.Item(14, 9).Value = 740.8333
$ws.Cells.Item(14, 11).Value = 2222.4999
$ws.Cells.Item(14, 13).Value = -2049.4999
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 386.875
$ws.Cells.Item(86, 9).Value = 386.875
$ws.Cells.Item(86, 11).Value = 1160.625
$ws.Cells.Item(86, 13).Value = 25.375
$ws.Cells.Item(89, 8).Value = 386.875
$ws.Cells.Item(89, 9).Value = 386.875
$ws.Cells.Item(89, 11).Value = 3481.875
$ws.Cells.Item(89, 13).Value = 2446.125
$ws.Cells.Item(117, 8).Value = 245.8
$ws.Cells.Item(117, 10).Value = 246.26315
$ws.Cells.Item(117, 12).Value = 738.78945
$ws.Cells.Item(117, 14).Value = -7622.78945
$ws.Cells.Item(119, 8).Value = 6757.25
$ws.Cells.Item(119, 9).Value = 3514.5
$ws.Cells.Item(119, 10).Value = 10000
$ws.Cells.Item(119, 11).Value = 10543.5
$ws.Cells.Item(119, 12).Value = 30000
$ws.Cells.Item(119, 13).Value = -5705.5
$ws.Cells.Item(119, 14).Value = -39676

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3617.818
$ws.Cells.Item(70, 9).Value = 3558.111
$ws.Cells.Item(70, 11).Value = 3558.111
$ws.Cells.Item(70, 13).Value = -3288.111
$ws.Cells.Item(73, 8).Value = 3617.818
$ws.Cells.Item(73, 9).Value = 3558.111
$ws.Cells.Item(73, 11).Value = 3558.111
$ws.Cells.Item(73, 13).Value = -2622.111

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 79168.664
$ws.Cells.Item(2, 9).Value = 49998.332
$ws.Cells.Item(2, 10).Value = 108339
$ws.Cells.Item(2, 11).Value = 49998.332
$ws.Cells.Item(2, 12).Value = 108339
$ws.Cells.Item(2, 13).Value = -49886.332
$ws.Cells.Item(2, 14).Value = -108563
$ws.Cells.Item(16, 8).Value = 2333
$ws.Cells.Item(16, 9).Value = 1999.5
$ws.Cells.Item(16, 10).Value = 3000
$ws.Cells.Item(16, 11).Value = 1999.5
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = -1829.5
$ws.Cells.Item(16, 14).Value = -3340
$ws.Cells.Item(46, 8).Value = 5732.222
$ws.Cells.Item(46, 9).Value = 2150
$ws.Cells.Item(46, 11).Value = 2150
$ws.Cells.Item(46, 13).Value = -1962
$ws.Cells.Item(55, 8).Value = 1020.38464
$ws.Cells.Item(55, 9).Value = 1016.5
$ws.Cells.Item(55, 10).Value = 1033.3334
$ws.Cells.Item(55, 11).Value = 1016.5
$ws.Cells.Item(55, 12).Value = 1033.3334
$ws.Cells.Item(55, 13).Value = -843.5
$ws.Cells.Item(55, 14).Value = -1379.3334
$ws.Cells.Item(127, 8).Value = 48833.332
$ws.Cells.Item(127, 10).Value = 48833.332
$ws.Cells.Item(127, 12).Value = 48833.332
$ws.Cells.Item(127, 14).Value = -58753.332
$ws.Cells.Item(136, 8).Value = 1525
$ws.Cells.Item(136, 9).Value = 1525
$ws.Cells.Item(136, 11).Value = 4575
$ws.Cells.Item(136, 13).Value = -2025

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 35217.43
$ws.Cells.Item(41, 10).Value = 35058.4
$ws.Cells.Item(41, 12).Value = 35058.4
$ws.Cells.Item(41, 14).Value = -35838.4
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).Value = $null
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 3830.0715
$ws.Cells.Item(122, 9).Value = 2933.1
$ws.Cells.Item(122, 10).Value = 6072.5
$ws.Cells.Item(122, 11).Value = 8799.299999999999
$ws.Cells.Item(122, 12).Value = 18217.5
$ws.Cells.Item(122, 13).Value = -6349.299999999999
$ws.Cells.Item(122, 14).Value = -23117.5
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = $null
$ws.Cells.Item(130, 8).Value = 40000
$ws.Cells.Item(130, 10).Value = 40000
$ws.Cells.Item(130, 12).Value = 40000
$ws.Cells.Item(130, 14).Value = -50040
$ws.Cells.Item(135, 8).Value = 38500
$ws.Cells.Item(135, 10).Value = 38500
$ws.Cells.Item(135, 12).Value = 38500
$ws.Cells.Item(135, 14).Value = -48640
